$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update D3/E3 and D4/E4 from TC01_... filenames to TC10_... filenames,
# matching the values already used in D2/E2.
$ws.Range("D3").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_WebData.xlsx"
$ws.Range("D4").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC10_Canine_Filter_SamplePatho-TCellLymphoma_WebData.xlsx"

# Update the active selection to D4:F4 with active cell D4.
$ws.Activate()
$ws.Range("D4:F4").Select()
